$d = $word.ActiveDocument

# --- Paragraph 2: "olvidado" passage ---
# The three runs making up this paragraph already concatenate to the
# correct text; re-enter the text as a whole so Word collapses the run
# boundaries into a single run (reapplying the original formatting).
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2NoMark = $d.Range($r2.Start, $r2.End - 1)
$p2FontName = $r2NoMark.Font.Name
$p2FontSize = $r2NoMark.Font.Size
$p2Text = $r2NoMark.Text
$r2NoMark.Delete()
$p2InsertAt = $d.Range($r2.Start, $r2.Start)
$p2InsertAt.InsertAfter($p2Text)

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2NoMark = $d.Range($r2.Start, $r2.End - 1)
$r2NoMark.Font.Name = $p2FontName
$r2NoMark.Font.Size = $p2FontSize

# --- Last paragraph: add a trailing date sentence ---
# "Los enanos se enteraron ... nunca regresaron." gains a trailing space
# and a brand-new run with the date sentence.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$lastFontName = $rLast.Font.Name
$lastFontSize = $rLast.Font.Size
$spaceInsertAt = $d.Range($rLast.End - 1, $rLast.End - 1)
$spaceInsertAt.InsertAfter(" ")

$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$paraBreakAt = $d.Range($rLast.End - 1, $rLast.End - 1)
$paraBreakAt.InsertParagraphAfter()

$dateSentence = "Estos hechos ocurrieron alrededor del año 350."
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParaRange = $newPara.Range
$newTextAt = $d.Range($newParaRange.Start, $newParaRange.Start)
$newTextAt.InsertAfter($dateSentence)

# Merge the newly created paragraph back into the previous one so the
# sentence becomes a second run within the same paragraph.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$rLast = $pLast.Range
$markRange = $d.Range($rLast.End - 1, $rLast.End)
$markRange.Delete()

# Make sure the new run keeps the same font as the rest of the paragraph.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$newRunRange = $d.Range($rLast.End - 1 - $dateSentence.Length, $rLast.End - 1)
$newRunRange.Font.Name = $lastFontName
$newRunRange.Font.Size = $lastFontSize
